$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the next day's results as a new row (row 41)
$ws.Range("A41").NumberFormat = $ws.Range("A40").NumberFormat
$ws.Range("A41").Value = 45990
$ws.Range("B41").Value = 91
$ws.Range("C41").Value = 102
$ws.Range("D41").Value = 98
